$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Subscript three (U+2083) used in PEPE price notation (row 35)
$sub3 = [char]0x2083

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.308.17"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.69%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.071.24"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.65%  "

# Row 4
$ws.Range("E4").Value = "  -0.24%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.53%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.62%  "

# Row 7
$ws.Range("E7").Value = "  -0.18%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.553"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.29%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.069.99"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.32%  "

# Row 10
$ws.Range("E10").Value = "  -3.21%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.88"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.14%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.462"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.41%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000241"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.47%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.16"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.41%  "

# Row 15
$ws.Range("E15").Value = "  -1.81%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.580.24"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.63%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.22"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.41%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.351.19"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.28%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.072.56"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.54%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "477.32"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.07%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.61"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.36%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.718"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.18%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.52"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.02%  "

# Row 24
$ws.Range("B24").Value = "Fetch.AI"
$ws.Range("C24").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.36"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.87%  "

# Row 25
$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.07"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.21%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "81.68"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.39%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.14%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.75"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.12%  "

# Row 29
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.67"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.80%  "

# Row 30
$ws.Range("B30").Value = "FirstDigitalUSD"
$ws.Range("C30").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.22%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.28"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.44%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.20"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.31%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.114"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.95%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "27.31"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.72%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0${sub3}0852"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.32%  "

# Row 36
$ws.Range("E36").Value = "  -1.88%  "

# Row 37
$ws.Range("B37").Value = "dogwifhat"
$ws.Range("C37").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.36"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.56%  "

# Row 38
$ws.Range("B38").Value = "Filecoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.12"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.92%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.21"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.37%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "9.38"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.42%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "50.24"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.63%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "443.86"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.08%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.286"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.00%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0363"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.23%  "

# Row 45
$ws.Range("B45").Value = "Arweave"
$ws.Range("C45").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "40.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.50%  "

# Row 46
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.806.89"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.18%  "

# Row 47
$ws.Range("E47").Value = "  +1.03%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "130.79"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.06%  "

# Row 49
$ws.Range("E49").Value = "  +0.03%  "

# Row 50
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.07"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.12%  "

# Row 51
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.112"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.98%  "
